# Adding test case for Search Module OPQA-1238
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("A94").Value = "TestCase_B93"
$ws.Range("B94").Value = "OPQA-1238"
$ws.Range("C94").Value = "Verify that profile page of a person gets displayed when clicks on any PEOPLE search result in ALL search results page"
$ws.Range("D94").Value = "Y"
$ws.Range("E94").Value = "PASS"

$ws.Range("A94:E94").Style = $ws.Range("A93:E93").Style

$ws.Range("D91").Select()
